# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-12 11:15:25
#
# The "Recorded By" column (G) contains a comma-separated list of the
# users/accounts that recorded each attendance session. In this sync, the
# relative position of the "System" entry within each multi-value list was
# swapped with its neighboring entry (this is how the upstream export now
# orders those tokens). Single-value cells, and multi-value cells that do
# not contain a "System" entry, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column G ("Recorded By").
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Text

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) {
        continue
    }

    $idx = [Array]::IndexOf($parts, "System")
    if ($idx -eq -1) {
        continue
    }

    if ($idx -eq ($parts.Count - 1)) {
        $other = $idx - 1
    } else {
        $other = $idx + 1
    }

    $tmp = $parts[$idx]
    $parts[$idx] = $parts[$other]
    $parts[$other] = $tmp

    $cell.Value = ($parts -join ", ")
}
